$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 machine number
$ws.Range("E2").Value = "213-GKD"

# Add new rows for Kerr 6th floor machines
$data = @(
    @(6, 2, "214-GKD"),
    @(6, 3, "215-GKD"),
    @(6, 4, "216-GKD")
)

$row = 3
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = "KER"
    $ws.Cells.Item($row, 2).Value = $entry[0]
    $ws.Cells.Item($row, 3).Value = $entry[1]
    $ws.Cells.Item($row, 4).Value = "A"
    $ws.Cells.Item($row, 5).Value = $entry[2]
    $row++
}

$ws.Range("E6").Select()
